# Applies the "Fruta, Femacal de La Calera - Uva" weekly update.
# Rows 591-596 get updated values (new week of data pushed in),
# and 4 brand-new rows are appended (597-600), shifting the former
# last row (597) down to 601. Final used range becomes A1:T601.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $D, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T) {
    $ws.Cells.Item($Row, 1).Value = 3
    $ws.Cells.Item($Row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($Row, 3).Value = "Coquimbo"

    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Range($ws.Cells.Item($Row, 4), $ws.Cells.Item($Row, 4)).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($Row, 5).Value = 5
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100109
    $ws.Cells.Item($Row, 8).Value = "Uva"
    $ws.Cells.Item($Row, 9).Value = 100109001
    $ws.Cells.Item($Row, 10).Value = "Uva"
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
}

# Existing rows 591-596: values refreshed with the new week's report.
Set-Row 591 44628 "Red Globe"         "Primera" 78 12000 12000 12000 "$/caja 12 kilos"    "Provincia de San Felipe de Aconcagua" 1000 12
Set-Row 592 44628 "Red Globe"         "Primera" 56 14000 14000 14000 "$/caja 18 kilos"    "Provincia de San Felipe de Aconcagua" 778  18
Set-Row 593 44628 "Ribier"            "Primera" 60 13000 13000 13000 "$/caja 15 kilos"    "Provincia de San Felipe de Aconcagua" 867  15
Set-Row 594 44628 "Thompson seedless" "Primera" 70 13000 13000 13000 "$/caja 15 kilos"    "Provincia de San Felipe de Aconcagua" 867  15
Set-Row 595 44335 "Crimpson Seedless" "Primera" 65 13000 13000 13000 "$/caja 15 kilos"    "Provincia de San Felipe de Aconcagua" 867  15
Set-Row 596 44335 "Red Globe"         "Primera" 75 11000 11000 11000 "$/caja 15 kilos"    "Provincia de San Felipe de Aconcagua" 733  15

# New rows 597-600, inserted ahead of the former last row.
Set-Row 597 44552 "Red Globe"         "Primera" 70 15000 15000 15000 "$/bandeja 10 kilos" "Provincia de Copiapó" 1500 10
Set-Row 598 44552 "Superior Seedless" "Primera" 68 15000 15000 15000 "$/bandeja 10 kilos" "Provincia de Copiapó" 1500 10
Set-Row 599 44544 "Red Globe"         "Primera" 50 20000 20000 20000 "$/bandeja 10 kilos" "Provincia de Copiapó" 2000 10
Set-Row 600 44544 "Superior Seedless" "Primera" 48 20000 20000 20000 "$/bandeja 10 kilos" "Provincia de Copiapó" 2000 10

# Former row 597, now shifted down to row 601 (content unchanged).
Set-Row 601 44160 "Flame Seedless"    "Primera" 50 20000 20000 20000 "$/bandeja 10 kilos" "Provincia de Copiapó" 2000 10
